# Reorder the "nearest station" columns so that columns L:M hold the
# "Estación más cercana 6/7" (and their header/data), pushing the existing
# "Inicio estación más cercana 1-5" columns (L:P) two columns to the right
# (into N:R). Columns S:T ("Inicio ... 6/7") are left untouched.
#
# Before:  L=Inicio1 M=Inicio2 N=Inicio3 O=Inicio4 P=Inicio5 Q=Estacion6 R=Estacion7
# After:   L=Estacion6 M=Estacion7 N=Inicio1 O=Inicio2 P=Inicio3 Q=Inicio4 R=Inicio5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 41

for ($r = 1; $r -le $lastRow; $r++) {
    $oldL = $ws.Cells.Item($r, 12).Text   # L - Inicio estación más cercana 1
    $oldM = $ws.Cells.Item($r, 13).Text   # M - Inicio estación más cercana 2
    $oldN = $ws.Cells.Item($r, 14).Text   # N - Inicio estación más cercana 3
    $oldO = $ws.Cells.Item($r, 15).Text   # O - Inicio estación más cercana 4
    $oldP = $ws.Cells.Item($r, 16).Text   # P - Inicio estación más cercana 5
    $oldQ = $ws.Cells.Item($r, 17).Text   # Q - Estación más cercana 6
    $oldR = $ws.Cells.Item($r, 18).Text   # R - Estación más cercana 7

    $ws.Cells.Item($r, 12).Value = $oldQ  # L <- Estación más cercana 6
    $ws.Cells.Item($r, 13).Value = $oldR  # M <- Estación más cercana 7
    $ws.Cells.Item($r, 14).Value = $oldL  # N <- Inicio estación más cercana 1
    $ws.Cells.Item($r, 15).Value = $oldM  # O <- Inicio estación más cercana 2
    $ws.Cells.Item($r, 16).Value = $oldN  # P <- Inicio estación más cercana 3
    $ws.Cells.Item($r, 17).Value = $oldO  # Q <- Inicio estación más cercana 4
    $ws.Cells.Item($r, 18).Value = $oldP  # R <- Inicio estación más cercana 5
}
